$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 8990
$ws.Range("J51").Value = 8990
$ws.Range("L51").Value = 8990
$ws.Range("N51").Value = -9958
$ws.Range("H53").Value = 505.125
$ws.Range("J53").Value = 773
$ws.Range("L53").Value = 773
$ws.Range("N53").Value = -2047
$ws.Range("H70").Value = 0
$ws.Range("I70").Value = 0
$ws.Range("K70").Value = 0
$ws.Range("M70").ClearContents()
$ws.Range("H73").Value = 0
$ws.Range("I73").Value = 0
$ws.Range("K73").Value = 0
$ws.Range("M73").ClearContents()
$ws.Range("H80").Value = 302.42856
$ws.Range("I80").Value = 190.2
$ws.Range("J80").Value = 583
$ws.Range("K80").Value = 570.5999999999999
$ws.Range("L80").Value = 1749
$ws.Range("M80").Value = 427.4000000000001
$ws.Range("N80").Value = -3745
$ws.Range("H83").Value = 302.42856
$ws.Range("I83").Value = 190.2
$ws.Range("J83").Value = 583
$ws.Range("K83").Value = 1711.8
$ws.Range("L83").Value = 5247
$ws.Range("M83").Value = 3280.2
$ws.Range("N83").Value = -15231
$ws.Range("H121").Value = 1866.2858
$ws.Range("J121").Value = 1866.2858
$ws.Range("L121").Value = 5598.857400000001
$ws.Range("N121").Value = -9092.857400000001
$ws.Range("H129").Value = 2495
$ws.Range("I129").Value = 2495
$ws.Range("J129").Value = 0
$ws.Range("K129").Value = 7485
$ws.Range("L129").Value = 0
$ws.Range("M129").Value = -2485
$ws.Range("N129").ClearContents()
$ws.Range("H132").Value = 1396.3334
$ws.Range("I132").Value = 1442.1428
$ws.Range("K132").Value = 4326.428400000001
$ws.Range("M132").Value = -1796.428400000001
$ws.Range("H137").Value = 825.1667
$ws.Range("I137").Value = 967.6667
$ws.Range("J137").Value = 682.6667
$ws.Range("K137").Value = 2903.0001
$ws.Range("L137").Value = 2048.0001
$ws.Range("M137").Value = -353.0001000000002
$ws.Range("N137").Value = -7148.0001
$ws.Range("H141").Value = 6997
$ws.Range("J141").Value = 6992.3335
$ws.Range("L141").Value = 20977.0005
$ws.Range("N141").Value = -31337.0005

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2258.8
$ws.Range("I45").Value = 2215.3333
$ws.Range("K45").Value = 2215.3333
$ws.Range("M45").Value = -1838.3333
$ws.Range("H122").Value = 2133
$ws.Range("I122").Value = 2499.5
$ws.Range("J122").Value = 1400
$ws.Range("K122").Value = 7498.5
$ws.Range("L122").Value = 4200
$ws.Range("M122").Value = -5048.5
$ws.Range("N122").Value = -9100

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2473.9375
$ws.Range("I86").Value = 2605.5334
$ws.Range("K86").Value = 2605.5334
$ws.Range("M86").Value = -1482.5334
$ws.Range("H89").Value = 2473.9375
$ws.Range("I89").Value = 2605.5334
$ws.Range("K89").Value = 13027.667
$ws.Range("M89").Value = -7411.666999999999
$ws.Range("H105").Value = 1266.25
$ws.Range("I105").Value = 1304.5714
$ws.Range("K105").Value = 1304.5714
$ws.Range("M105").Value = 442.4286

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 1611.1538
$ws.Range("I122").Value = 1636.091
$ws.Range("K122").Value = 4908.272999999999
$ws.Range("M122").Value = -2458.272999999999

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 49122230
$ws.Range("I4").Value = 4250931
$ws.Range("J4").Value = 89007830
$ws.Range("K4").Value = 12752793
$ws.Range("L4").Value = 267023490
$ws.Range("M4").Value = -12752681
$ws.Range("N4").Value = -267023714
$ws.Range("H26").Value = 1003.2778
$ws.Range("I26").Value = 186.33333
$ws.Range("J26").Value = 1166.6666
$ws.Range("K26").Value = 558.99999
$ws.Range("L26").Value = 3499.9998
$ws.Range("M26").Value = -270.99999
$ws.Range("N26").Value = -4075.9998
$ws.Range("H131").Value = 2666.6667
$ws.Range("J131").Value = 3500
$ws.Range("L131").Value = 10500
$ws.Range("N131").Value = -20580
$ws.Range("H137").Value = 849.6667
$ws.Range("I137").Value = 849.6667
$ws.Range("J137").Value = 0
$ws.Range("K137").Value = 2549.0001
$ws.Range("L137").Value = 0
$ws.Range("M137").Value = 2550.9999
$ws.Range("N137").ClearContents()
$ws.Range("H140").Value = 2775.4285
$ws.Range("I140").Value = 2654.8333
$ws.Range("J140").Value = 3499
$ws.Range("K140").Value = 7964.499899999999
$ws.Range("L140").Value = 10497
$ws.Range("M140").Value = -2784.499899999999
$ws.Range("N140").Value = -20857
$ws.Range("H141").Value = 700
$ws.Range("I141").Value = 700
$ws.Range("K141").Value = 2100
$ws.Range("M141").Value = 3080

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H32").Value = 0
$ws.Range("J32").Value = 0
$ws.Range("L32").Value = 0
$ws.Range("N32").ClearContents()
$ws.Range("H80").Value = 9316.666999999999
$ws.Range("I80").Value = 9180
$ws.Range("J80").Value = 10000
$ws.Range("K80").Value = 9180
$ws.Range("L80").Value = 10000
$ws.Range("M80").Value = -8182
$ws.Range("N80").Value = -11996
$ws.Range("H83").Value = 9316.666999999999
$ws.Range("I83").Value = 9180
$ws.Range("J83").Value = 10000
$ws.Range("K83").Value = 45900
$ws.Range("L83").Value = 50000
$ws.Range("M83").Value = -40908
$ws.Range("N83").Value = -59984
$ws.Range("H132").Value = 3525.3
$ws.Range("I132").Value = 2894.3333
$ws.Range("K132").Value = 8682.999899999999
$ws.Range("M132").Value = -6152.999899999999

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5000
$ws.Range("I7").Value = 5000
$ws.Range("K7").Value = 5000
$ws.Range("M7").Value = -4888
$ws.Range("H16").Value = 4997.7144
$ws.Range("I16").Value = 4197
$ws.Range("K16").Value = 4197
$ws.Range("M16").Value = -4027
$ws.Range("H46").Value = 0
$ws.Range("I46").Value = 0
$ws.Range("K46").Value = 0
$ws.Range("M46").ClearContents()
$ws.Range("H82").Value = 2299.1667
$ws.Range("J82").Value = 2531.6667
$ws.Range("L82").Value = 2531.6667
$ws.Range("N82").Value = -3253.6667
$ws.Range("H85").Value = 2299.1667
$ws.Range("J85").Value = 2531.6667
$ws.Range("L85").Value = 2531.6667
$ws.Range("N85").Value = -5027.6667
$ws.Range("H93").Value = 617
$ws.Range("I93").Value = 617
$ws.Range("J93").Value = 0
$ws.Range("K93").Value = 617
$ws.Range("L93").Value = 0
$ws.Range("M93").Value = 631
$ws.Range("N93").ClearContents()
$ws.Range("H126").Value = 5000
$ws.Range("I126").Value = 5000
$ws.Range("K126").Value = 15000
$ws.Range("M126").Value = -12530
$ws.Range("H131").Value = 0
$ws.Range("J131").Value = 0
$ws.Range("L131").Value = 0
$ws.Range("N131").ClearContents()

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 13333
$ws.Range("J2").Value = 16000
$ws.Range("L2").Value = 16000
$ws.Range("N2").Value = -16224
$ws.Range("H100").Value = 1450
$ws.Range("I100").Value = 1450
$ws.Range("K100").Value = 2900
$ws.Range("M100").Value = -2359
$ws.Range("H136").Value = 499.66666
$ws.Range("I136").Value = 499.66666
$ws.Range("K136").Value = 1498.99998
$ws.Range("M136").Value = 1051.00002
